# Apply updated loading_percent values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.719734031184064
$ws.Range("C2").Value = 6.338110475625068
$ws.Range("D2").Value = 5.997737432143495
$ws.Range("E2").Value = 16.3666063605383
$ws.Range("G2").Value = 42.38067966786666
$ws.Range("H2").Value = 17.13778305690583
$ws.Range("I2").Value = 25.16676969106781
$ws.Range("K2").Value = 10.26130066688681
$ws.Range("B3").Value = 9.461937575043839
$ws.Range("C3").Value = 6.050277162415481
$ws.Range("D3").Value = 5.881881094603054
$ws.Range("E3").Value = 15.44511530997996
$ws.Range("G3").Value = 41.85692823163191
$ws.Range("H3").Value = 17.11128556180114
$ws.Range("I3").Value = 25.09955508138161
$ws.Range("K3").Value = 10.07635003867921
$ws.Range("B4").Value = 9.303349880089117
$ws.Range("C4").Value = 5.868497479506255
$ws.Range("D4").Value = 5.811461299113171
$ws.Range("E4").Value = 14.8556570611034
$ws.Range("G4").Value = 41.54362538032703
$ws.Range("H4").Value = 17.09815159638737
$ws.Range("I4").Value = 25.06303788995626
$ws.Range("K4").Value = 9.96442590369794
$ws.Range("B5").Value = 9.238763880175611
$ws.Range("C5").Value = 5.793274114758037
$ws.Range("D5").Value = 5.782987892451319
$ws.Range("E5").Value = 14.60977272385425
$ws.Range("G5").Value = 41.41816597367848
$ws.Range("H5").Value = 17.09358864705591
$ws.Range("I5").Value = 25.04935580087666
$ws.Range("K5").Value = 9.919296313628495
$ws.Range("B6").Value = 9.228045278371811
$ws.Range("C6").Value = 5.78071798562221
$ws.Range("D6").Value = 5.778274649123396
$ws.Range("E6").Value = 14.56861009125208
$ws.Range("G6").Value = 41.39747098839726
$ws.Range("H6").Value = 17.09287865280647
$ws.Range("I6").Value = 25.04715640319368
$ws.Range("K6").Value = 9.911833604648832
$ws.Range("B7").Value = 9.302478534258773
$ws.Range("C7").Value = 5.86748745944276
$ws.Range("D7").Value = 5.811076335889142
$ws.Range("E7").Value = 14.85236355709098
$ws.Range("G7").Value = 41.54192425704269
$ws.Range("H7").Value = 17.09808686273232
$ws.Range("I7").Value = 25.06284850993723
$ws.Range("K7").Value = 9.963815234247269
$ws.Range("B8").Value = 9.630975589730228
$ws.Range("C8").Value = 6.239982276833132
$ws.Range("D8").Value = 5.957667225279609
$ws.Range("E8").Value = 16.05393537092442
$ws.Range("G8").Value = 42.19845145791776
$ws.Range("H8").Value = 17.12799655310679
$ws.Range("I8").Value = 25.14261027148558
$ws.Range("K8").Value = 10.19723145888315
$ws.Range("B9").Value = 10.26798162196165
$ws.Range("C9").Value = 6.925766823442977
$ws.Range("D9").Value = 6.248944812509253
$ws.Range("E9").Value = 18.25653853013034
$ws.Range("G9").Value = 43.54562022148843
$ws.Range("H9").Value = 17.2114733983449
$ws.Range("I9").Value = 25.33651921295178
$ws.Range("K9").Value = 10.66492333023957
$ws.Range("B10").Value = 10.72561598882068
$ws.Range("C10").Value = 7.397188808125008
$ws.Range("D10").Value = 6.462825292422861
$ws.Range("E10").Value = 19.88668104026993
$ws.Range("G10").Value = 44.56352712681013
$ws.Range("H10").Value = 17.28783155285596
$ws.Range("I10").Value = 25.50148649833267
$ws.Range("K10").Value = 11.01063080371898
$ws.Range("B11").Value = 10.93041924948407
$ws.Range("C11").Value = 7.603765392600851
$ws.Range("D11").Value = 6.559606875506117
$ws.Range("E11").Value = 20.58681388573462
$ws.Range("G11").Value = 45.03091777679401
$ws.Range("H11").Value = 17.32579761288815
$ws.Range("I11").Value = 25.58132095522448
$ws.Range("K11").Value = 11.1675475373749
$ws.Range("B12").Value = 11.00740267345625
$ws.Range("C12").Value = 7.68080182586435
$ws.Range("D12").Value = 6.596143028984638
$ws.Range("E12").Value = 20.84601655382136
$ws.Range("G12").Value = 45.20838027984539
$ws.Range("H12").Value = 17.34063492466881
$ws.Range("I12").Value = 25.61223060731812
$ws.Range("K12").Value = 11.22685375078801
$ws.Range("B13").Value = 10.99084960089599
$ws.Range("C13").Value = 7.664264374034051
$ws.Range("D13").Value = 6.588279937588309
$ws.Range("E13").Value = 20.79045547298146
$ws.Range("G13").Value = 45.17014197757977
$ws.Range("H13").Value = 17.33741904003912
$ws.Range("I13").Value = 25.60554368519461
$ws.Range("K13").Value = 11.21408719992606
$ws.Range("B14").Value = 10.93676463226868
$ws.Range("C14").Value = 7.610127363197858
$ws.Range("D14").Value = 6.562615206674945
$ws.Range("E14").Value = 20.60825714161747
$ws.Range("G14").Value = 45.04550913437576
$ws.Range("H14").Value = 17.32700908628997
$ws.Range("I14").Value = 25.58385038649717
$ws.Range("K14").Value = 11.17242935551052
$ws.Range("B15").Value = 10.90355923613941
$ws.Range("C15").Value = 7.576810416194507
$ws.Range("D15").Value = 6.546878986194431
$ws.Range("E15").Value = 20.49588515298019
$ws.Range("G15").Value = 44.9692248705509
$ws.Range("H15").Value = 17.32069251281655
$ws.Range("I15").Value = 25.57065058354055
$ws.Range("K15").Value = 11.14689590699502
$ws.Range("B16").Value = 10.71215618366535
$ws.Range("C16").Value = 7.383524968274828
$ws.Range("D16").Value = 6.456486666452445
$ws.Range("E16").Value = 19.84009506517228
$ws.Range("G16").Value = 44.53305632909887
$ws.Range("H16").Value = 17.28541504442726
$ws.Range("I16").Value = 25.49636462432182
$ws.Range("K16").Value = 11.00036333039968
$ws.Range("B17").Value = 10.59380852524149
$ws.Range("C17").Value = 7.262889821357295
$ws.Range("D17").Value = 6.400873910768332
$ws.Range("E17").Value = 19.42720337214358
$ws.Range("G17").Value = 44.26648383711328
$ws.Range("H17").Value = 17.26459807432017
$ws.Range("I17").Value = 25.45201214871608
$ws.Range("K17").Value = 10.91033435338694
$ws.Range("B18").Value = 10.52542488042343
$ws.Range("C18").Value = 7.19276516188804
$ws.Range("D18").Value = 6.368840786295917
$ws.Range("E18").Value = 19.18581760584146
$ws.Range("G18").Value = 44.11357592045924
$ws.Range("H18").Value = 17.25292904417333
$ws.Range("I18").Value = 25.42695334865465
$ws.Range("K18").Value = 10.85852307613733
$ws.Range("B19").Value = 10.50222028654717
$ws.Range("C19").Value = 7.168897211009694
$ws.Range("D19").Value = 6.357988219198734
$ws.Range("E19").Value = 19.1034174053626
$ws.Range("G19").Value = 44.06188029862092
$ws.Range("H19").Value = 17.249030498226
$ws.Range("I19").Value = 25.41854673603532
$ws.Range("K19").Value = 10.8409777110257
$ws.Range("B20").Value = 10.60643989357185
$ws.Range("C20").Value = 7.275808515214487
$ws.Range("D20").Value = 6.40679904797121
$ws.Range("E20").Value = 19.47155992309255
$ws.Range("G20").Value = 44.294818828826
$ws.Range("H20").Value = 17.26678260653501
$ws.Range("I20").Value = 25.45668688406138
$ws.Range("K20").Value = 10.91992157099377
$ws.Range("B21").Value = 10.95266684420996
$ws.Range("C21").Value = 7.626061412008545
$ws.Range("D21").Value = 6.570156926426961
$ws.Range("E21").Value = 20.66193366576245
$ws.Range("G21").Value = 45.0821052222431
$ws.Range("H21").Value = 17.3300542829904
$ws.Range("I21").Value = 25.59020392216459
$ws.Range("K21").Value = 11.18466889004744
$ws.Range("B22").Value = 11.17558266976858
$ws.Range("C22").Value = 7.848017522988543
$ws.Range("D22").Value = 6.676245834407564
$ws.Range("E22").Value = 21.40540715340152
$ws.Range("G22").Value = 45.59932315677695
$ws.Range("H22").Value = 17.37408716179697
$ws.Range("I22").Value = 25.68141164847878
$ws.Range("K22").Value = 11.35699983972689
$ws.Range("B23").Value = 11.05694243328055
$ws.Range("C23").Value = 7.730208264249273
$ws.Range("D23").Value = 6.619698315309669
$ws.Range("E23").Value = 21.01174710722376
$ws.Range("G23").Value = 45.32307875968179
$ws.Range("H23").Value = 17.35034216972701
$ws.Range("I23").Value = 25.63237506956685
$ws.Range("K23").Value = 11.26510794167082
$ws.Range("B24").Value = 10.60073031683528
$ws.Range("C24").Value = 7.269970367398872
$ws.Range("D24").Value = 6.404120479262844
$ws.Range("E24").Value = 19.45151881612521
$ws.Range("G24").Value = 44.282007482904
$ws.Range("H24").Value = 17.26579404803737
$ws.Range("I24").Value = 25.45457206469583
$ws.Range("K24").Value = 10.91558734820281
$ws.Range("B25").Value = 10.09706787748221
$ws.Range("C25").Value = 6.745613349238916
$ws.Range("D25").Value = 6.169995430674483
$ws.Range("E25").Value = 17.65112784128387
$ws.Range("G25").Value = 43.17562627590854
$ws.Range("H25").Value = 17.18624168711098
$ws.Range("I25").Value = 25.28007960135419
$ws.Range("K25").Value = 10.53775034334967
